$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 6.240107999999999
$ws.Cells.Item(2, 8).Value = 18.720324
$ws.Cells.Item(2, 9).Value = 0.01732230523539376
$ws.Cells.Item(2, 10).Value = 0.01732230523539376
$ws.Cells.Item(2, 13).Value = 0.3883076666666667
$ws.Cells.Item(2, 14).Value = 1.164923
$ws.Cells.Item(2, 15).Value = 0.1188638477168776
$ws.Cells.Item(2, 16).Value = 0.1188638477168776
$ws.Cells.Item(2, 17).Value = 2.423081777228
$ws.Cells.Item(2, 18).Value = 21.807735995052
$ws.Cells.Item(2, 19).Value = 0.002058995851605116
$ws.Cells.Item(2, 20).Value = 0.002058995851605116

# Row 3
$ws.Cells.Item(3, 7).Value = 6.240107999999999
$ws.Cells.Item(3, 8).Value = 18.720324
$ws.Cells.Item(3, 9).Value = 0.01732230523539376
$ws.Cells.Item(3, 10).Value = 0.01732230523539376
$ws.Cells.Item(3, 15).Value = 0.6829215134520935
$ws.Cells.Item(3, 16).Value = 0.6829215134520935
$ws.Cells.Item(3, 17).Value = 13.92159774656
$ws.Cells.Item(3, 18).Value = 125.29437971904
$ws.Cells.Item(3, 19).Value = 0.01182977490783423
$ws.Cells.Item(3, 20).Value = 0.01182977490783423

# Row 4
$ws.Cells.Item(4, 7).Value = 6.240107999999999
$ws.Cells.Item(4, 8).Value = 18.720324
$ws.Cells.Item(4, 9).Value = 0.01732230523539376
$ws.Cells.Item(4, 10).Value = 0.01732230523539376
$ws.Cells.Item(4, 15).Value = 0.1982146388310289
$ws.Cells.Item(4, 16).Value = 0.1982146388310289
$ws.Cells.Item(4, 17).Value = 4.040675853563999
$ws.Cells.Item(4, 18).Value = 36.366082682076
$ws.Cells.Item(4, 19).Value = 0.003433534475954416
$ws.Cells.Item(4, 20).Value = 0.003433534475954416

# Row 5
$ws.Cells.Item(5, 9).Value = 0.9592798330716089
$ws.Cells.Item(5, 10).Value = 0.9592798330716091
$ws.Cells.Item(5, 13).Value = 0.3883076666666667
$ws.Cells.Item(5, 14).Value = 1.164923
$ws.Cells.Item(5, 15).Value = 0.1188638477168776
$ws.Cells.Item(5, 16).Value = 0.1188638477168776
$ws.Cells.Item(5, 17).Value = 134.1861519694723
$ws.Cells.Item(5, 18).Value = 1207.675367725251
$ws.Cells.Item(5, 19).Value = 0.1140236919960955
$ws.Cells.Item(5, 20).Value = 0.1140236919960955

# Row 6
$ws.Cells.Item(6, 9).Value = 0.9592798330716089
$ws.Cells.Item(6, 10).Value = 0.9592798330716091
$ws.Cells.Item(6, 15).Value = 0.6829215134520935
$ws.Cells.Item(6, 16).Value = 0.6829215134520935
$ws.Cells.Item(6, 19).Value = 0.6551128354253348
$ws.Cells.Item(6, 20).Value = 0.6551128354253349

# Row 7
$ws.Cells.Item(7, 9).Value = 0.9592798330716089
$ws.Cells.Item(7, 10).Value = 0.9592798330716091
$ws.Cells.Item(7, 15).Value = 0.1982146388310289
$ws.Cells.Item(7, 16).Value = 0.1982146388310289
$ws.Cells.Item(7, 19).Value = 0.1901433056501787
$ws.Cells.Item(7, 20).Value = 0.1901433056501787

# Row 8
$ws.Cells.Item(8, 7).Value = 8.428738666666666
$ws.Cells.Item(8, 9).Value = 0.02339786169299727
$ws.Cells.Item(8, 10).Value = 0.02339786169299728
$ws.Cells.Item(8, 13).Value = 0.3883076666666667
$ws.Cells.Item(8, 14).Value = 1.164923
$ws.Cells.Item(8, 15).Value = 0.1188638477168776
$ws.Cells.Item(8, 16).Value = 0.1188638477168776
$ws.Cells.Item(8, 17).Value = 3.272943844596444
$ws.Cells.Item(8, 18).Value = 29.456494601368
$ws.Cells.Item(8, 19).Value = 0.002781159869176992
$ws.Cells.Item(8, 20).Value = 0.002781159869176992

# Row 9
$ws.Cells.Item(9, 7).Value = 8.428738666666666
$ws.Cells.Item(9, 9).Value = 0.02339786169299727
$ws.Cells.Item(9, 10).Value = 0.02339786169299728
$ws.Cells.Item(9, 15).Value = 0.6829215134520935
$ws.Cells.Item(9, 16).Value = 0.6829215134520935
$ws.Cells.Item(9, 19).Value = 0.01597890311892446
$ws.Cells.Item(9, 20).Value = 0.01597890311892446

# Row 10
$ws.Cells.Item(10, 7).Value = 8.428738666666666
$ws.Cells.Item(10, 9).Value = 0.02339786169299727
$ws.Cells.Item(10, 10).Value = 0.02339786169299728
$ws.Cells.Item(10, 15).Value = 0.1982146388310289
$ws.Cells.Item(10, 16).Value = 0.1982146388310289
$ws.Cells.Item(10, 17).Value = 5.457886435042666
$ws.Cells.Item(10, 19).Value = 0.004637798704895821
$ws.Cells.Item(10, 20).Value = 0.004637798704895823

